$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new "period" data row (2509) below the existing last data row (2508, row 23) ---
# Insert a new blank row at 24, shifting everything below (blank rows + signature lines) down by one.
$ws.Rows.Item(24).Insert()

# Copy row 23 (the last existing data row, "2508") - values + formatting - into the new row 24.
# This mirrors Excel's behaviour when a user duplicates the last row to add a new period.
$ws.Range("B23:J23").Copy($ws.Range("B24:J24"))

# The row that used to be "last" (23) is no longer the last row, so give it the normal
# "interior" row formatting (matching rows 16-22) instead of the special bottom-border style.
$ws.Range("B22:J22").Copy()
$ws.Range("B23:J23").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Set the new period value for the freshly inserted row.
$ws.Cells.Item(24, 5).Value2 = "2509"

# --- Update summary figures to reflect the new period ---
# Cant. Periodos: 8 -> 9
$ws.Cells.Item(13, 6).Value2 = 9

# Valor Mora (total): 416000 -> 468000 (52000 * 9 periods)
$ws.Cells.Item(11, 5).Value2 = 468000
